# Refresh the crypto price/volume/coin table to the latest scrape.
#
# Column D ("Price") holds numeric-looking strings (e.g. "1.004",
# "24.789.81") that must stay plain TEXT, exactly as scraped, rather than
# being auto-coerced into numbers (which would drop trailing zeros / merge
# the multi-dot thousand-separator values). We force text by prefixing a
# leading apostrophe (the same quote-prefix Excel itself applies when a
# user types `'1.004` into a cell), built up via string concatenation so
# the PowerShell quoting stays simple.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.789.81'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '1.701.05'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").Value = "'" + '1.004'
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").Value = "'" + '316.56'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").Value = "'" + '1.007'
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("D8").Value = "'" + '0.4031'
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").Value = "'" + '1.504'
$ws.Range("E9").Value = '  -2.63%  '
$ws.Range("D10").Value = "'" + '54.03'
$ws.Range("E10").Value = '  -1.82%  '
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = "'" + '0.08892'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = "'" + '7.222'
$ws.Range("E13").Value = '  -1.40%  '
$ws.Range("D14").Value = "'" + '23.38'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = "'" + '8.018'
$ws.Range("E15").Value = '  +4.78%  '
$ws.Range("D16").Value = "'" + '0.00001325'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("D17").Value = '1.715.80'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = "'" + '100.00'
$ws.Range("D19").Value = "'" + '0.07036'
$ws.Range("E19").Value = '  -0.26%  '
$ws.Range("D20").Value = "'" + '19.64'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = "'" + '7.027'
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").Value = "'" + '1.003'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("E23").Value = '  +2.42%  '
$ws.Range("D24").Value = '24.751.15'
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").Value = "'" + '3.201'
$ws.Range("E25").Value = '  +7.11%  '
$ws.Range("D26").Value = "'" + '2.357'
$ws.Range("E26").Value = '  +1.55%  '
$ws.Range("D27").Value = "'" + '22.76'
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").Value = "'" + '162.09'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").Value = "'" + '136.34'
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("D30").Value = "'" + '5.171'
$ws.Range("D31").Value = "'" + '7.805'
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").Value = "'" + '0.08744'
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").Value = "'" + '1.072'
$ws.Range("E33").Value = '  -3.59%  '
$ws.Range("D34").Value = "'" + '7.168'
$ws.Range("E34").Value = '  -3.72%  '
$ws.Range("D35").Value = "'" + '11.20'
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("D36").Value = "'" + '1.979'
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("D37").Value = "'" + '0.2738'
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").Value = "'" + '14.37'
$ws.Range("E38").Value = '  -3.17%  '
$ws.Range("D39").Value = "'" + '0.09187'
$ws.Range("E39").Value = '  +1.61%  '
$ws.Range("D40").Value = "'" + '0.02744'
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("D41").Value = "'" + '1.460'
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("D42").Value = "'" + '0.7665'
$ws.Range("E42").Value = '  -1.60%  '
$ws.Range("B43").Value = 'Decentraland'
$ws.Range("C43").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D43").Value = "'" + '0.7148'
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'" + '15.67'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").Value = "'" + '2.570'
$ws.Range("E45").Value = '  +1.98%  '
$ws.Range("D46").Value = "'" + '4.210'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("D47").Value = "'" + '1.004'
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").Value = "'" + '140.58'
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").Value = "'" + '1.313'
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("D50").Value = "'" + '90.85'
$ws.Range("E50").Value = '  +2.60%  '
$ws.Range("D51").Value = "'" + '0.07985'
$ws.Range("E51").Value = '  -0.64%  '
